# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the AHB-diff header columns from the generic "_old"/"_new" suffixes
# to the concrete format-version suffixes ("_FV2210"/"_FV2304"), wraps the
# sheet's used range in an Excel Table (ListObject) and freezes the header
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row -------------------------------------------------
$oldSuffix = "_old"
$newSuffix = "_new"
$fv2210 = "_FV2210"
$fv2304 = "_FV2304"

$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i] + $fv2210
}

# Column 11 ("diff") is unchanged.

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $headers[$i] + $fv2304
}

# --- 2. Wrap data range in a Table ----------------------------------------
$rng = $ws.Range("A1:U61")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
